# Auto-generated: update Price (D) and Volume(1h) (E) columns for crypto rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$c = $ws.Range("D2")
$s = $c.Style
$c.Value = "'26.283.33"
$c.Style = $s
$c = $ws.Range("E2")
$s = $c.Style
$c.Value = "  +2.96%  "
$c.Style = $s

# Row 3: Ethereum
$c = $ws.Range("D3")
$s = $c.Style
$c.Value = "'1.720.95"
$c.Style = $s
$c = $ws.Range("E3")
$s = $c.Style
$c.Value = "  +3.31%  "
$c.Style = $s

# Row 4: TetherUSD
$c = $ws.Range("D4")
$s = $c.Style
$c.Value = "'0.9991"
$c.Style = $s
$c = $ws.Range("E4")
$s = $c.Style
$c.Value = "  -0.01%  "
$c.Style = $s

# Row 5: BNB
$c = $ws.Range("D5")
$s = $c.Style
$c.Value = "'240.17"
$c.Style = $s
$c = $ws.Range("E5")
$s = $c.Style
$c.Value = "  +1.08%  "
$c.Style = $s

# Row 6: USDC
$c = $ws.Range("D6")
$s = $c.Style
$c.Value = "'0.9997"
$c.Style = $s
$c = $ws.Range("E6")
$s = $c.Style
$c.Value = "  -0.07%  "
$c.Style = $s

# Row 7: XRP
$c = $ws.Range("D7")
$s = $c.Style
$c.Value = "'0.4732"
$c.Style = $s
$c = $ws.Range("E7")
$s = $c.Style
$c.Value = "  -1.32%  "
$c.Style = $s

# Row 8: Cardano
$c = $ws.Range("D8")
$s = $c.Style
$c.Value = "'0.2628"
$c.Style = $s
$c = $ws.Range("E8")
$s = $c.Style
$c.Value = "  +0.07%  "
$c.Style = $s

# Row 9: Dogecoin
$c = $ws.Range("E9")
$s = $c.Style
$c.Value = "  +0.50%  "
$c.Style = $s

# Row 10: WrappedEther
$c = $ws.Range("D10")
$s = $c.Style
$c.Value = "'1.716.73"
$c.Style = $s
$c = $ws.Range("E10")
$s = $c.Style
$c.Value = "  +3.13%  "
$c.Style = $s

# Row 11: TRON
$c = $ws.Range("D11")
$s = $c.Style
$c.Value = "'0.07052"
$c.Style = $s
$c = $ws.Range("E11")
$s = $c.Style
$c.Value = "  -0.47%  "
$c.Style = $s

# Row 12: Solana
$c = $ws.Range("D12")
$s = $c.Style
$c.Value = "'15.51"
$c.Style = $s
$c = $ws.Range("E12")
$s = $c.Style
$c.Value = "  +4.76%  "
$c.Style = $s

# Row 13: Polygon
$c = $ws.Range("D13")
$s = $c.Style
$c.Value = "'0.5985"
$c.Style = $s
$c = $ws.Range("E13")
$s = $c.Style
$c.Value = "  +1.93%  "
$c.Style = $s

# Row 14: Polkadot
$c = $ws.Range("D14")
$s = $c.Style
$c.Value = "'4.433"
$c.Style = $s
$c = $ws.Range("E14")
$s = $c.Style
$c.Value = "  +1.61%  "
$c.Style = $s

# Row 15: Litecoin
$c = $ws.Range("D15")
$s = $c.Style
$c.Value = "'76.26"
$c.Style = $s
$c = $ws.Range("E15")
$s = $c.Style
$c.Value = "  +1.74%  "
$c.Style = $s

# Row 16: Dai
$c = $ws.Range("D16")
$s = $c.Style
$c.Value = "'0.9997"
$c.Style = $s
$c = $ws.Range("E16")
$s = $c.Style
$c.Value = "  -0.03%  "
$c.Style = $s

# Row 17: BinanceUSD
$c = $ws.Range("E17")
$s = $c.Style
$c.Value = "  -0.01%  "
$c.Style = $s

# Row 18: WrappedBTC
$c = $ws.Range("D18")
$s = $c.Style
$c.Value = "'26.292.11"
$c.Style = $s
$c = $ws.Range("E18")
$s = $c.Style
$c.Value = "  +3.00%  "
$c.Style = $s

# Row 19: ShibaInu
$c = $ws.Range("D19")
$s = $c.Style
$c.Value = "'0.000006810"
$c.Style = $s
$c = $ws.Range("E19")
$s = $c.Style
$c.Value = "  +0.88%  "
$c.Style = $s

# Row 20: Avalanche
$c = $ws.Range("D20")
$s = $c.Style
$c.Value = "'11.54"
$c.Style = $s
$c = $ws.Range("E20")
$s = $c.Style
$c.Value = "  +0.77%  "
$c.Style = $s

# Row 21: WrappedliquidstakedEther2.0
$c = $ws.Range("D21")
$s = $c.Style
$c.Value = "'1.936.14"
$c.Style = $s
$c = $ws.Range("E21")
$s = $c.Style
$c.Value = "  +3.13%  "
$c.Style = $s

# Row 22: Uniswap
$c = $ws.Range("E22")
$s = $c.Style
$c.Value = "  +2.78%  "
$c.Style = $s

# Row 23: Cosmos
$c = $ws.Range("D23")
$s = $c.Style
$c.Value = "'8.717"
$c.Style = $s
$c = $ws.Range("E23")
$s = $c.Style
$c.Value = "  -0.08%  "
$c.Style = $s

# Row 24: Chainlink
$c = $ws.Range("D24")
$s = $c.Style
$c.Value = "'5.244"
$c.Style = $s
$c = $ws.Range("E24")
$s = $c.Style
$c.Value = "  -0.58%  "
$c.Style = $s

# Row 25: Monero
$c = $ws.Range("D25")
$s = $c.Style
$c.Value = "'135.18"
$c.Style = $s
$c = $ws.Range("E25")
$s = $c.Style
$c.Value = "  -0.35%  "
$c.Style = $s

# Row 26: EthereumClassic
$c = $ws.Range("D26")
$s = $c.Style
$c.Value = "'15.21"
$c.Style = $s
$c = $ws.Range("E26")
$s = $c.Style
$c.Value = "  +1.28%  "
$c.Style = $s

# Row 27: LidoDAOToken
$c = $ws.Range("D27")
$s = $c.Style
$c.Value = "'1.767"
$c.Style = $s
$c = $ws.Range("E27")
$s = $c.Style
$c.Value = "  +3.43%  "
$c.Style = $s

# Row 28: Toncoin
$c = $ws.Range("D28")
$s = $c.Style
$c.Value = "'1.399"
$c.Style = $s
$c = $ws.Range("E28")
$s = $c.Style
$c.Value = "  +0.77%  "
$c.Style = $s

# Row 29: BitcoinCash
$c = $ws.Range("D29")
$s = $c.Style
$c.Value = "'106.64"
$c.Style = $s
$c = $ws.Range("E29")
$s = $c.Style
$c.Value = "  +1.75%  "
$c.Style = $s

# Row 30: InternetComputer(DFINITY)
$c = $ws.Range("D30")
$s = $c.Style
$c.Value = "'3.946"
$c.Style = $s
$c = $ws.Range("E30")
$s = $c.Style
$c.Value = "  -0.54%  "
$c.Style = $s

# Row 31: Filecoin
$c = $ws.Range("D31")
$s = $c.Style
$c.Value = "'3.688"
$c.Style = $s
$c = $ws.Range("E31")
$s = $c.Style
$c.Value = "  +1.41%  "
$c.Style = $s

# Row 32: Stellar
$c = $ws.Range("D32")
$s = $c.Style
$c.Value = "'0.07805"
$c.Style = $s
$c = $ws.Range("E32")
$s = $c.Style
$c.Value = "  +0.59%  "
$c.Style = $s

# Row 33: Hedera
$c = $ws.Range("D33")
$s = $c.Style
$c.Value = "'0.04507"
$c.Style = $s
$c = $ws.Range("E33")
$s = $c.Style
$c.Value = "  +7.05%  "
$c.Style = $s

# Row 34: HuobiToken
$c = $ws.Range("D34")
$s = $c.Style
$c.Value = "'2.613"
$c.Style = $s
$c = $ws.Range("E34")
$s = $c.Style
$c.Value = "  +0.53%  "
$c.Style = $s

# Row 35: ARBITRUM
$c = $ws.Range("D35")
$s = $c.Style
$c.Value = "'0.9811"
$c.Style = $s
$c = $ws.Range("E35")
$s = $c.Style
$c.Value = "  +3.62%  "
$c.Style = $s

# Row 36: ImmutableX
$c = $ws.Range("D36")
$s = $c.Style
$c.Value = "'0.6222"
$c.Style = $s
$c = $ws.Range("E36")
$s = $c.Style
$c.Value = "  +2.10%  "
$c.Style = $s

# Row 37: TrustWalletToken
$c = $ws.Range("D37")
$s = $c.Style
$c.Value = "'0.9340"
$c.Style = $s
$c = $ws.Range("E37")
$s = $c.Style
$c.Value = "  +8.35%  "
$c.Style = $s

# Row 38: Quant
$c = $ws.Range("D38")
$s = $c.Style
$c.Value = "'114.86"
$c.Style = $s
$c = $ws.Range("E38")
$s = $c.Style
$c.Value = "  +18.20%  "
$c.Style = $s

# Row 39: MXToken
$c = $ws.Range("D39")
$s = $c.Style
$c.Value = "'2.451"
$c.Style = $s
$c = $ws.Range("E39")
$s = $c.Style
$c.Value = "  -5.52%  "
$c.Style = $s

# Row 40: RenderToken
$c = $ws.Range("D40")
$s = $c.Style
$c.Value = "'1.927"
$c.Style = $s
$c = $ws.Range("E40")
$s = $c.Style
$c.Value = "  +4.28%  "
$c.Style = $s

# Row 41: PaxDollar
$c = $ws.Range("D41")
$s = $c.Style
$c.Value = "'0.9998"
$c.Style = $s
$c = $ws.Range("E41")
$s = $c.Style
$c.Value = "  -0.09%  "
$c.Style = $s

# Row 42: FraxShare
$c = $ws.Range("D42")
$s = $c.Style
$c.Value = "'5.658"
$c.Style = $s
$c = $ws.Range("E42")
$s = $c.Style
$c.Value = "  +17.02%  "
$c.Style = $s

# Row 43: VeChain
$c = $ws.Range("E43")
$s = $c.Style
$c.Value = "  +1.41%  "
$c.Style = $s

# Row 44: TheSandbox
$c = $ws.Range("D44")
$s = $c.Style
$c.Value = "'0.3834"
$c.Style = $s
$c = $ws.Range("E44")
$s = $c.Style
$c.Value = "  +2.12%  "
$c.Style = $s

# Row 45: Algorand
$c = $ws.Range("D45")
$s = $c.Style
$c.Value = "'0.1184"
$c.Style = $s
$c = $ws.Range("E45")
$s = $c.Style
$c.Value = "  +5.70%  "
$c.Style = $s

# Row 46: Aptos
$c = $ws.Range("D46")
$s = $c.Style
$c.Value = "'6.367"
$c.Style = $s
$c = $ws.Range("E46")
$s = $c.Style
$c.Value = "  +2.79%  "
$c.Style = $s

# Row 47: Cronos
$c = $ws.Range("D47")
$s = $c.Style
$c.Value = "'0.05266"
$c.Style = $s
$c = $ws.Range("E47")
$s = $c.Style
$c.Value = "  +0.03%  "
$c.Style = $s

# Row 48: EnergySwap
$c = $ws.Range("D48")
$s = $c.Style
$c.Value = "'7.816"
$c.Style = $s
$c = $ws.Range("E48")
$s = $c.Style
$c.Value = "  +6.89%  "
$c.Style = $s

# Row 49: Elrond
$c = $ws.Range("D49")
$s = $c.Style
$c.Value = "'30.40"
$c.Style = $s
$c = $ws.Range("E49")
$s = $c.Style
$c.Value = "  +2.31%  "
$c.Style = $s

# Row 50: Decentraland
$c = $ws.Range("D50")
$s = $c.Style
$c.Value = "'0.3386"
$c.Style = $s
$c = $ws.Range("E50")
$s = $c.Style
$c.Value = "  +1.87%  "
$c.Style = $s

# Row 51: NEARProtocol
$c = $ws.Range("D51")
$s = $c.Style
$c.Value = "'1.218"
$c.Style = $s
$c = $ws.Range("E51")
$s = $c.Style
$c.Value = "  +2.11%  "
$c.Style = $s
